# Obhajoba: Tabulka s rozptylenym polom
#
# Slide 11's title reads (as two runs joined by a soft line break):
#   "Algoritmy" <soft line break> "Tabuľka s rozptýlenými polom"
# The trailing word "polom" is corrected to "položkami", i.e. the title
# becomes:
#   "Algoritmy" <soft line break> "Tabuľka s rozptýlenými položkami"
#
# Because the replacement text differs from the word it replaces,
# PowerPoint keeps the untouched leading text ("Tabuľka s rozptýlenými ")
# in the original run and puts the edited word in a new run that carries
# the same run properties - exactly the two-run split shown in the diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)

$titleShape = $s.Shapes.Title
$titleRange = $titleShape.TextFrame.TextRange

$oldWord = "polom"
$newWord = "položkami"

$fullText = $titleRange.Text
$startIndex = $fullText.IndexOf($oldWord)

# 1-based character position expected by TextRange.Characters(Start, Length)
$wordRange = $titleRange.Characters($startIndex + 1, $oldWord.Length)
$wordRange.Text = $newWord
